# Slide 3, "TextBox 5" (How do we take notes? -> bullet list):
#   - split the second bullet's run so "information" becomes its own run
#     (matches what PowerPoint does when the user re-enters that run and
#     the proofing pass re-tags it)
#   - append a new third bullet paragraph "Blah "
#   - the textbox uses spAutoFit / wrap="none", so its box grows to fit
#     the new line; nudge Width/Height to the exact resulting extent

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(3)
$sh = $s.Shapes.Item(4)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# --- split "...'correct' information" into two runs ------------------
# Full text is 107 chars: "Take best notes, without any mistakes" (37)
# + CR + "Take notes, with mistakes and only transfer the 'correct'
# information" (69). "information" is the trailing 11 characters
# (1-based start 97).
$infoStart = $tr.Length - "information".Length + 1
$infoRange = $tr.Characters($infoStart, "information".Length)
$infoRange.Text = "information"

# --- add the new bullet paragraph -------------------------------------
$null = $tr.InsertAfter("`rBlah ")

# --- match the shape's resulting size exactly -------------------------
# (spAutoFit height differs from the target by a hair, and the width
# needs to grow slightly as well; set both explicitly in points)
$sh.Width = 540.3408203125
$sh.Height = 72.70315551757812
$null = 0
